$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Thanh_Toan" (Payment) description table lives in E18:G21 and had a
# row describing "t_tt (Tiền khách thanh toán)" - the amount the customer
# pays in - at E19:G19. Remove that row (no more customer-entered balance)
# and shift the remaining rows of that table up by one.
$ws.Range("E19:G20").Value2 = $ws.Range("E20:G21").Value2
$ws.Range("E21:G21").Clear()

# I32 had a one-off style (bold-less red font, fontId 10) that was only
# used by that single cell; align it with the neighbouring J32 cell's
# style so the now-unused style entry drops out of use.
$ws.Range("I32").Font.Color = $ws.Range("J32").Font.Color

# The active selection/scroll position changed as part of the edit.
$ws.Range("F22").Select()
